# Move example values to the correct columns on the "3ASY04_Genomics" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3ASY04_Genomics")

# Row 2: library preparation kit version was blank but should hold "v2";
# the value "HiSeq 2000" had been placed in the wrong column (ngs platform)
# and the library strand value "R1" was missing from its column.
$ws.Range("T2").Value = "v2"
$ws.Range("AC2").Value = ""
$ws.Range("AO2").Value = "R1"

# Row 3: "R2" had been entered under base-calling software version instead
# of library strand.
$ws.Range("AL3").Value = ""
$ws.Range("AO3").Value = "R2"

# Row 4: "Forward" had been entered under base-calling software version
# instead of library strand.
$ws.Range("AL4").Value = ""
$ws.Range("AO4").Value = "Forward"
